# Applies the price/volume/coin-order refresh described by the commit
# "Updated cryptos list on Thu Aug 29 21:33:18 UTC 2024 with GitHub Actions".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '''59.395.80'
$ws.Range('E2').Value = '  +0.37%  '
# Row 3
$ws.Range('D3').Value = '''2.525.66'
$ws.Range('E3').Value = '  +0.52%  '
# Row 4
$ws.Range('E4').Value = '  +0.05%  '
# Row 5
$ws.Range('D5').Value = '''537.33'
$ws.Range('E5').Value = '  -0.14%  '
# Row 6
$ws.Range('D6').Value = '''139.74'
$ws.Range('E6').Value = '  -3.14%  '
# Row 7
$ws.Range('E7').Value = '  +0.35%  '
# Row 8
$ws.Range('D8').Value = '''0.563'
$ws.Range('E8').Value = '  -1.58%  '
# Row 9
$ws.Range('D9').Value = '''2.529.66'
$ws.Range('E9').Value = '  -0.79%  '
# Row 10
$ws.Range('D10').Value = '''0.0992'
$ws.Range('E10').Value = '  -0.78%  '
# Row 11
$ws.Range('E11').Value = '  +1.34%  '
# Row 12
$ws.Range('E12').Value = '  -3.06%  '
# Row 13
$ws.Range('E13').Value = '  +0.81%  '
# Row 14
$ws.Range('D14').Value = '''2.972.81'
$ws.Range('E14').Value = '  +0.76%  '
# Row 15
$ws.Range('D15').Value = '''23.12'
$ws.Range('E15').Value = '  -2.56%  '
# Row 16
$ws.Range('D16').Value = '''59.325.88'
$ws.Range('E16').Value = '  +0.47%  '
# Row 17
$ws.Range('E17').Value = '  +0.67%  '
# Row 18
$ws.Range('D18').Value = '''2.522.22'
$ws.Range('E18').Value = '  -0.42%  '
# Row 19
$ws.Range('D19').Value = '''10.94'
$ws.Range('E19').Value = '  -3.63%  '
# Row 20
$ws.Range('E20').Value = '  -1.85%  '
# Row 21
$ws.Range('D21').Value = '''321.69'
$ws.Range('E21').Value = '  -0.48%  '
# Row 22
$ws.Range('D22').Value = '''0.998'
$ws.Range('E22').Value = '  +0.14%  '
# Row 23
$ws.Range('D23').Value = '''5.79'
$ws.Range('E23').Value = '  +0.20%  '
# Row 24
$ws.Range('D24').Value = '''61.29'
$ws.Range('E24').Value = '  -1.26%  '
# Row 25
$ws.Range('E25').Value = '  -4.37%  '
# Row 26
$ws.Range('E26').Value = '  +1.41%  '
# Row 27
$ws.Range('D27').Value = '''0.998'
$ws.Range('E27').Value = '  +0.75%  '
# Row 28
$ws.Range('E28').Value = '  +0.10%  '
# Row 29
$ws.Range('D29').Value = '''6.74'
$ws.Range('E29').Value = '  -0.78%  '
# Row 30
$ws.Range('B30').Value = 'PEPE'
$ws.Range('C30').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D30').Value = '''0.0₃0769'
$ws.Range('E30').Value = '  -1.27%  '
# Row 31
$ws.Range('B31').Value = 'PancakeSwap'
$ws.Range('C31').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D31').Value = '''1.80'
$ws.Range('E31').Value = '  +0.49%  '
# Row 32
$ws.Range('D32').Value = '''160.32'
$ws.Range('E32').Value = '  +0.87%  '
# Row 33
$ws.Range('E33').Value = '  +0.37%  '
# Row 34
$ws.Range('B34').Value = 'ImmutableX'
$ws.Range('C34').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D34').Value = '''1.46'
$ws.Range('E34').Value = '  +1.04%  '
# Row 35
$ws.Range('B35').Value = 'Fetch.AI'
$ws.Range('C35').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D35').Value = '''1.13'
$ws.Range('E35').Value = '  -6.97%  '
# Row 36
$ws.Range('E36').Value = '  -0.39%  '
# Row 37
$ws.Range('E37').Value = '  -4.61%  '
# Row 38
$ws.Range('D38').Value = '''1.59'
$ws.Range('E38').Value = '  -1.87%  '
# Row 39
$ws.Range('D39').Value = '''37.00'
$ws.Range('E39').Value = '  +0.34%  '
# Row 40
$ws.Range('E40').Value = '  -0.52%  '
# Row 41
$ws.Range('D41').Value = '''0.811'
$ws.Range('E41').Value = '  -1.65%  '
# Row 42
$ws.Range('B42').Value = 'Bittensor'
$ws.Range('C42').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D42').Value = '''284.75'
$ws.Range('E42').Value = '  -6.83%  '
# Row 43
$ws.Range('B43').Value = 'RenderToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D43').Value = '''5.26'
$ws.Range('E43').Value = '  -8.13%  '
# Row 44
$ws.Range('E44').Value = '  +0.46%  '
# Row 45
$ws.Range('D45').Value = '''0.599'
$ws.Range('E45').Value = '  -0.82%  '
# Row 46
$ws.Range('D46').Value = '''10.87'
$ws.Range('E46').Value = '  +0.84%  '
# Row 47
$ws.Range('D47').Value = '''123.75'
$ws.Range('E47').Value = '  -1.54%  '
# Row 48
$ws.Range('E48').Value = '  -0.92%  '
# Row 49
$ws.Range('D49').Value = '''18.54'
$ws.Range('E49').Value = '  -1.22%  '
# Row 50
$ws.Range('E50').Value = '  -1.86%  '
# Row 51
$ws.Range('D51').Value = '''0.0222'
$ws.Range('E51').Value = '  -2.46%  '
